$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the "Toollink" (D) column for the rows that previously had no
# link text (rows 6-12) with the same "Null" placeholder used in D4/D5.
for ($r = 6; $r -le 12; $r++) {
    $ws.Range("D$r").Value = "Null"
}

# Seed D3 then D2 with their final text first so the shared-string table
# picks up "https://www.google.nl/" (index 24) before the combined/garbled
# "https://Google.nl/,https://stackoverflow.com/..." text (index 25),
# matching the order the strings were authored in.
$ws.Range("D3").Value = "https://www.google.nl/"
$ws.Range("D2").Value = "https://Google.nl/,https://stackoverflow.com/questions/179713/how-to-change-the-href-attribute-for-a-hyperlink-using-jquery"

# Now turn D2 and D3 into real hyperlinks (D2 first so it claims rId5,
# then D3 claims rId6), and re-apply the workbook's existing "Hyperlink"
# cell style so the cells match the look of the other link cells.
$ws.Hyperlinks.Add($ws.Range("D2"), "https://Google.nl/,https://stackoverflow.com/questions/179713/how-to-change-the-href-attribute-for-a-hyperlink-using-jquery")
$ws.Range("D2").Style = "Hyperlink"

$ws.Hyperlinks.Add($ws.Range("D3"), "https://www.google.nl/")
$ws.Range("D3").Style = "Hyperlink"

# Update the active selection to reflect where work left off.
$ws.Range("B13").Select()
